$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (the current "id" column),
# shifting id/name/capacity/material_nature/kwargs one column to the right.
$ws.Range("B1").EntireColumn.Insert()

# Set the new header cell value and copy the style from the neighboring header cell.
$ws.Range("B1").Value = "env"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
